# Update duplicate-detection comparison values for rows 2 and 3 (duplistid pair)
# These are the randomly re-sampled "noise" comparison metrics/flags used by
# the duplicates-check dofile (sesion 5, semana 3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 1
$ws.Cells.Item(2, 14).Value = 1
$ws.Cells.Item(2, 15).Value = 5.4117140769958496
$ws.Cells.Item(2, 16).Value = 41.535858154296875
$ws.Cells.Item(2, 17).Value = 15.914735794067383
$ws.Cells.Item(2, 18).Value = 25.621122360229492
$ws.Cells.Item(2, 19).Value = 10.464022636413574
$ws.Cells.Item(2, 20).Value = 37.738624572753906
$ws.Cells.Item(2, 21).Value = 9.9870309829711914
$ws.Cells.Item(2, 22).Value = 27.751594543457031
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 1
$ws.Cells.Item(2, 25).Value = 6.1416792869567871
$ws.Cells.Item(2, 26).Value = 30.753311157226563
$ws.Cells.Item(2, 27).Value = 30.753311157226563
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = 17.944118499755859
$ws.Cells.Item(2, 30).Value = 29.970684051513672
$ws.Cells.Item(2, 31).Value = 5.6441683769226074
$ws.Cells.Item(2, 32).Value = 24.326515197753906
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 35).Value = 10.537778854370117
$ws.Cells.Item(2, 36).Value = 12.444417953491211
$ws.Cells.Item(2, 37).Value = 7.4625377655029297
$ws.Cells.Item(2, 38).Value = 4.9818801879882813
$ws.Cells.Item(2, 39).Value = 19.420539855957031
$ws.Cells.Item(2, 40).Value = 34.3587646484375
$ws.Cells.Item(2, 41).Value = 23.82172966003418
$ws.Cells.Item(2, 42).Value = 10.53703498840332
$ws.Cells.Item(2, 44).Value = 1
$ws.Cells.Item(2, 45).Value = 17.804193496704102
$ws.Cells.Item(2, 46).Value = 2.6125538349151611
$ws.Cells.Item(2, 47).Value = 2.6125538349151611
$ws.Cells.Item(2, 48).Value = 0
$ws.Cells.Item(2, 49).Value = 2.83695387840271
$ws.Cells.Item(2, 50).Value = 31.2218017578125
$ws.Cells.Item(2, 51).Value = 14.377996444702148
$ws.Cells.Item(2, 52).Value = 16.843805313110352
$ws.Cells.Item(2, 55).Value = 4.3230438232421875
$ws.Cells.Item(2, 56).Value = 9.7958612442016602
$ws.Cells.Item(2, 57).Value = 9.7958612442016602
$ws.Cells.Item(2, 58).Value = 0
$ws.Cells.Item(2, 59).Value = 13.057753562927246
$ws.Cells.Item(2, 60).Value = 5.6126728057861328
$ws.Cells.Item(2, 61).Value = 5.6126728057861328
$ws.Cells.Item(2, 62).Value = 0
$ws.Cells.Item(2, 64).Value = 1
$ws.Cells.Item(2, 65).Value = 18.734275817871094
$ws.Cells.Item(2, 66).Value = 9.6759710311889648
$ws.Cells.Item(2, 67).Value = 9.6759710311889648
$ws.Cells.Item(2, 68).Value = 0
$ws.Cells.Item(2, 69).Value = 14.227289199829102
$ws.Cells.Item(2, 70).Value = 34.413803100585938
$ws.Cells.Item(2, 71).Value = 20.556005477905273
$ws.Cells.Item(2, 72).Value = 13.857797622680664
$ws.Cells.Item(2, 73).Value = 2.328934907913208
$ws.Cells.Item(2, 74).Value = 8.7675819396972656
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 1
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 15).Value = 5.4117140769958496
$ws.Cells.Item(3, 16).Value = 41.535858154296875
$ws.Cells.Item(3, 17).Value = 15.914735794067383
$ws.Cells.Item(3, 18).Value = 25.621122360229492
$ws.Cells.Item(3, 19).Value = 10.464022636413574
$ws.Cells.Item(3, 20).Value = 37.738624572753906
$ws.Cells.Item(3, 21).Value = 9.9870309829711914
$ws.Cells.Item(3, 22).Value = 27.751594543457031
$ws.Cells.Item(3, 23).Value = 0
$ws.Cells.Item(3, 24).Value = 1
$ws.Cells.Item(3, 25).Value = 6.1416792869567871
$ws.Cells.Item(3, 26).Value = 30.753311157226563
$ws.Cells.Item(3, 27).Value = 30.753311157226563
$ws.Cells.Item(3, 28).Value = 0
$ws.Cells.Item(3, 29).Value = 17.944118499755859
$ws.Cells.Item(3, 30).Value = 29.970684051513672
$ws.Cells.Item(3, 31).Value = 5.6441683769226074
$ws.Cells.Item(3, 32).Value = 24.326515197753906
$ws.Cells.Item(3, 33).Value = 0
$ws.Cells.Item(3, 35).Value = 10.537778854370117
$ws.Cells.Item(3, 36).Value = 12.444417953491211
$ws.Cells.Item(3, 37).Value = 7.4625377655029297
$ws.Cells.Item(3, 38).Value = 4.9818801879882813
$ws.Cells.Item(3, 39).Value = 19.420539855957031
$ws.Cells.Item(3, 40).Value = 34.3587646484375
$ws.Cells.Item(3, 41).Value = 23.82172966003418
$ws.Cells.Item(3, 42).Value = 10.53703498840332
$ws.Cells.Item(3, 44).Value = 1
$ws.Cells.Item(3, 45).Value = 17.804193496704102
$ws.Cells.Item(3, 46).Value = 2.6125538349151611
$ws.Cells.Item(3, 47).Value = 2.6125538349151611
$ws.Cells.Item(3, 48).Value = 0
$ws.Cells.Item(3, 49).Value = 2.83695387840271
$ws.Cells.Item(3, 50).Value = 31.2218017578125
$ws.Cells.Item(3, 51).Value = 14.377996444702148
$ws.Cells.Item(3, 52).Value = 16.843805313110352
$ws.Cells.Item(3, 55).Value = 4.3230438232421875
$ws.Cells.Item(3, 56).Value = 9.7958612442016602
$ws.Cells.Item(3, 57).Value = 9.7958612442016602
$ws.Cells.Item(3, 58).Value = 0
$ws.Cells.Item(3, 59).Value = 13.057753562927246
$ws.Cells.Item(3, 60).Value = 5.6126728057861328
$ws.Cells.Item(3, 61).Value = 5.6126728057861328
$ws.Cells.Item(3, 62).Value = 0
$ws.Cells.Item(3, 64).Value = 1
$ws.Cells.Item(3, 65).Value = 18.734275817871094
$ws.Cells.Item(3, 66).Value = 9.6759710311889648
$ws.Cells.Item(3, 67).Value = 9.6759710311889648
$ws.Cells.Item(3, 68).Value = 0
$ws.Cells.Item(3, 69).Value = 14.227289199829102
$ws.Cells.Item(3, 70).Value = 34.413803100585938
$ws.Cells.Item(3, 71).Value = 20.556005477905273
$ws.Cells.Item(3, 72).Value = 13.857797622680664
$ws.Cells.Item(3, 73).Value = 2.328934907913208
$ws.Cells.Item(3, 74).Value = 8.7675819396972656